# Added gerstung score and HSCT prediction
# Inserts a new gene row for "ATRX" into the "Genetics (core)" block,
# just above the existing "RUNX1" row, pushing every subsequent row
# down by one (old row 54 "RUNX1" becomes row 55, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 54 (shifts RUNX1.. down to 55..103)
$ws.Rows("54:54").Insert()

# Fill in the new ATRX row, mirroring the layout of the ASXL1 row above it
$ws.Range("A54").Value = $ws.Range("A53").Text
$ws.Range("B54").Value = "ATRX"
$ws.Range("C54").Value = "ATRX"
$ws.Range("D54").Value = $ws.Range("D53").Text
$ws.Range("E54").Value = $ws.Range("E53").Text
$ws.Range("F54").Value = $ws.Range("F53").Text
$ws.Range("G54").Value = $ws.Range("G53").Text

# Match the slightly shorter row height recorded for the new row
$ws.Rows("54:54").RowHeight = 13.8

# Reflect the cursor position left behind by the edit
$ws.Range("G54").Select()
